$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.541.89"
$ws.Range("E2").Value = "  -7.69%  "

$ws.Range("D3").Value = "2.550.44"
$ws.Range("E3").Value = "  -1.89%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'295.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.08%  "

$ws.Range("D6").Value = "'90.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.30%  "

$ws.Range("D7").Value = "'0.571"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.39%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.86%  "

$ws.Range("D10").Value = "'35.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.95%  "

$ws.Range("D11").Value = "'0.0801"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.40%  "

$ws.Range("D12").Value = "'7.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.41%  "

$ws.Range("D13").Value = "2.936.76"
$ws.Range("E13").Value = "  -1.97%  "

$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").Value = "2.560.05"
$ws.Range("E15").Value = "  -1.57%  "

$ws.Range("D16").Value = "'0.860"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.95%  "

$ws.Range("D17").Value = "'14.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.29%  "

$ws.Range("D18").Value = "42.541.26"
$ws.Range("E18").Value = "  -7.84%  "

$ws.Range("D19").Value = "'6.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("D20").Value = "0.0₃0965"
$ws.Range("E20").Value = "  -5.17%  "

$ws.Range("D21").Value = "'12.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.60%  "

$ws.Range("D22").Value = "'72.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").Value = "'257.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.47%  "

$ws.Range("D24").Value = "'2.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.79%  "

$ws.Range("D25").Value = "'29.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.41%  "

$ws.Range("D26").Value = "'2.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.22%  "

$ws.Range("E27").Value = "  +0.27%  "

$ws.Range("D28").Value = "'9.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.38%  "

$ws.Range("E29").Value = "  -4.40%  "

$ws.Range("D30").Value = "'35.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.14%  "

$ws.Range("D31").Value = "'5.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.52%  "

$ws.Range("D32").Value = "'150.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.28%  "

$ws.Range("D33").Value = "'2.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.71%  "

$ws.Range("D34").Value = "'3.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.69%  "

$ws.Range("D36").Value = "'0.0787"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.26%  "

$ws.Range("E37").Value = "  -7.81%  "

$ws.Range("D38").Value = "'24.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.06%  "

$ws.Range("E39").Value = "  -3.82%  "

$ws.Range("D40").Value = "'15.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.34%  "

$ws.Range("D41").Value = "'3.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.99%  "

$ws.Range("E42").Value = "  -6.87%  "

$ws.Range("D43").Value = "'3.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.36%  "

$ws.Range("D44").Value = "2.060.33"
$ws.Range("E44").Value = "  -1.63%  "

$ws.Range("D45").Value = "'0.998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").Value = "'84.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -12.19%  "

$ws.Range("E47").Value = "  +2.65%  "

$ws.Range("D48").Value = "2.792.87"
$ws.Range("E48").Value = "  -2.02%  "

$ws.Range("D49").Value = "'8.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.37%  "

$ws.Range("E50").Value = "  -3.44%  "

$ws.Range("D51").Value = "'102.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.21%  "
